$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Add the new result string for MaxPressure average travel time (D5)
$ws.Range("D5").Value = "46.6 saniye"

# Move selection to E5 to reflect the active editing location
$ws.Range("E5").Select()

# Mark the MaxPressure checkbox (Check Box 15 / ctrlProp1) as checked
try {
    $ws.Shapes("Check Box 15").ControlFormat.Value = 1
} catch {
}
